$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" '42.809.34'
Set-TextValue "E2" '  +0.16%  '
Set-TextValue "D3" '2.312.59'
Set-TextValue "E3" '  +0.80%  '
Set-TextValue "E4" '  +0.00%  '
Set-TextValue "D5" '301.85'
Set-TextValue "E5" '  -0.45%  '
Set-TextValue "D6" '95.24'
Set-TextValue "E6" '  -1.15%  '
Set-TextValue "D7" '0.505'
Set-TextValue "E7" '  +0.13%  '
Set-TextValue "E8" '  +0.03%  '
Set-TextValue "D9" '0.492'
Set-TextValue "D10" '34.20'
Set-TextValue "E10" '  -1.96%  '
Set-TextValue "D11" '18.90'
Set-TextValue "E11" '  +1.82%  '
Set-TextValue "E12" '  +0.14%  '
Set-TextValue "E13" '  -0.02%  '
Set-TextValue "D14" '6.73'
Set-TextValue "E14" '  -1.52%  '
Set-TextValue "D15" '2.675.86'
Set-TextValue "E15" '  +0.90%  '
Set-TextValue "D16" '2.315.25'
Set-TextValue "E16" '  +0.17%  '
Set-TextValue "D17" '0.786'
Set-TextValue "E17" '  +1.67%  '
Set-TextValue "D18" '42.757.74'
Set-TextValue "E18" '  +0.22%  '
Set-TextValue "D19" '12.15'
Set-TextValue "E19" '  -4.88%  '
Set-TextValue "D20" '6.13'
Set-TextValue "E20" '  +2.25%  '
Set-TextValue "D21" '0.0₃0890'
Set-TextValue "E21" '  -0.28%  '
Set-TextValue "D22" '67.74'
Set-TextValue "E22" '  +0.90%  '
Set-TextValue "E23" '  +6.29%  '
Set-TextValue "D24" '235.38'
Set-TextValue "E24" '  -0.21%  '
Set-TextValue "E25" '  +0.00%  '
Set-TextValue "D26" '2.41'
Set-TextValue "E26" '  +1.08%  '
Set-TextValue "D27" '24.32'
Set-TextValue "E27" '  -1.44%  '
Set-TextValue "E28" '  +14.80%  '
Set-TextValue "D29" '165.84'
Set-TextValue "E29" '  -1.02%  '
Set-TextValue "D30" '9.12'
Set-TextValue "E30" '  +1.41%  '
Set-TextValue "D31" '32.00'
Set-TextValue "E31" '  -2.79%  '
Set-TextValue "E32" '  -0.05%  '
Set-TextValue "D33" '4.99'
Set-TextValue "E33" '  +0.87%  '
Set-TextValue "D34" '17.71'
Set-TextValue "E34" '  +0.06%  '
Set-TextValue "D35" '4.45'
Set-TextValue "E35" '  +0.40%  '
Set-TextValue "D36" '0.0696'
Set-TextValue "E36" '  +1.77%  '
Set-TextValue "D37" '2.34'
Set-TextValue "E37" '  -0.88%  '
Set-TextValue "E38" '  +2.07%  '
Set-TextValue "E39" '  +0.13%  '
Set-TextValue "D40" '2.72'
Set-TextValue "E40" '  +1.14%  '
Set-TextValue "E41" '  -0.35%  '
Set-TextValue "D42" '20.97'
Set-TextValue "E42" '  +14.89%  '
Set-TextValue "D43" '1.926.33'
Set-TextValue "E43" '  -3.38%  '
Set-TextValue "E44" '  -0.29%  '
Set-TextValue "D45" '10.13'
Set-TextValue "E45" '  +0.34%  '
Set-TextValue "E46" '  -0.58%  '
Set-TextValue "D47" '2.74'
Set-TextValue "E47" '  -0.91%  '
Set-TextValue "D48" '2.88'
Set-TextValue "E48" '  -0.24%  '
Set-TextValue "D49" '2.544.23'
Set-TextValue "E49" '  +1.01%  '
Set-TextValue "D50" '53.17'
Set-TextValue "E50" '  -0.76%  '
Set-TextValue "D51" '72.02'
Set-TextValue "E51" '  +1.76%  '
